$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 105, shifting existing rows 105-168 down to 106-169.
$ws.Rows(105).Insert()

# Populate the newly inserted row 105 with a new weekly data point.
$ws.Cells.Item(105, 1).Value = 5
$ws.Cells.Item(105, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(105, 3).Value = "Maule"
$ws.Cells.Item(105, 4).Value = 45029
$ws.Cells.Item(105, 5).Value = 7
$ws.Cells.Item(105, 6).Value = "Fruta"
$ws.Cells.Item(105, 7).Value = 100108
$ws.Cells.Item(105, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(105, 9).Value = 100108002
$ws.Cells.Item(105, 10).Value = "Mango"
$ws.Cells.Item(105, 11).Value = "Sin especificar"
$ws.Cells.Item(105, 12).Value = "Primera"
$ws.Cells.Item(105, 13).Value = 248
$ws.Cells.Item(105, 14).Value = 7000
$ws.Cells.Item(105, 15).Value = 7000
$ws.Cells.Item(105, 16).Value = 7000
$ws.Cells.Item(105, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(105, 18).Value = "Perú"
$ws.Cells.Item(105, 19).Value = 1750
$ws.Cells.Item(105, 20).Value = 4

# Make sure the date column keeps the same custom date/time number format
# used by the rest of column D.
$ws.Cells.Item(105, 4).NumberFormat = $ws.Cells.Item(106, 4).NumberFormat
